$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, without altering its visible style,
# even when the text looks like a number (e.g. "44.83") or would otherwise be
# auto-converted by Excel's type inference.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# Plain text / URL / percent-string updates (safe to assign directly).
$ws.Range("D2").Value = '29.019.82'
$ws.Range("D3").Value = '1.829.33'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  -5.35%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E9").Value = '  +7.17%  '
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '1.829.80'
$ws.Range("E13").Value = '  -4.90%  '
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("E17").Value = '  +5.90%  '
$ws.Range("E18").Value = '  -2.12%  '
$ws.Range("D19").Value = '28.798.11'
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("E21").Value = '  -1.18%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -2.42%  '
$ws.Range("E27").Value = '  -2.85%  '
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("E34").Value = '  -1.38%  '
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").Value = '1.277.26'
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("E41").Value = '  +5.58%  '
$ws.Range("E42").Value = '  -3.23%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  -1.09%  '
$ws.Range("D45").Value = '1.979.81'
$ws.Range("E45").Value = '  -2.63%  '
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("E49").Value = '  -15.89%  '
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  +1.34%  '

# Numeric-looking price values that must remain stored as text.
Set-TextValue $ws.Range("D4") '0.9988'
Set-TextValue $ws.Range("D5") '241.09'
Set-TextValue $ws.Range("D6") '0.6239'
Set-TextValue $ws.Range("D8") '0.07557'
Set-TextValue $ws.Range("D9") '44.83'
Set-TextValue $ws.Range("D11") '22.73'
Set-TextValue $ws.Range("D12") '0.07641'
Set-TextValue $ws.Range("D14") '4.950'
Set-TextValue $ws.Range("D15") '0.6639'
Set-TextValue $ws.Range("D16") '82.20'
Set-TextValue $ws.Range("D17") '0.000009104'
Set-TextValue $ws.Range("D18") '5.984'
Set-TextValue $ws.Range("D20") '224.67'
Set-TextValue $ws.Range("D21") '12.32'
Set-TextValue $ws.Range("D23") '7.190'
Set-TextValue $ws.Range("D25") '159.80'
Set-TextValue $ws.Range("D26") '8.396'
Set-TextValue $ws.Range("D27") '0.1357'
Set-TextValue $ws.Range("D29") '1.494'
Set-TextValue $ws.Range("D30") '4.026'
Set-TextValue $ws.Range("D31") '4.042'
Set-TextValue $ws.Range("D32") '1.200'
Set-TextValue $ws.Range("D33") '0.05199'
Set-TextValue $ws.Range("D34") '1.842'
Set-TextValue $ws.Range("D35") '1.150'
Set-TextValue $ws.Range("D36") '0.7308'
Set-TextValue $ws.Range("D39") '2.756'
Set-TextValue $ws.Range("D40") '0.01787'
Set-TextValue $ws.Range("D41") '6.395'
Set-TextValue $ws.Range("D42") '0.8912'
Set-TextValue $ws.Range("D44") '101.57'
Set-TextValue $ws.Range("D46") '0.5106'
Set-TextValue $ws.Range("D47") '63.54'
Set-TextValue $ws.Range("D48") '0.00000000120'
Set-TextValue $ws.Range("D49") '0.07352'
Set-TextValue $ws.Range("D50") '0.3975'
Set-TextValue $ws.Range("D51") '8.873'
